$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix stray non-breaking space before newline in the "Clear / Visibility: 14 km" English string (C19)
$ws.Range("C19").Value = "[name=""""]   Clear \ Visibility: 14 km `n"

# Replace escaped double-quoted "Captain" with single-quoted 'Captain' in English column (C) dialogue lines
$ws.Range("C68").Value = "[name=""'Captain'""]   I'm not thirsty! I don't need your water!`n"
$ws.Range("C70").Value = "[name=""'Captain'""]   What kind of question is that? Did you forget why we're even out here in the first place?`n"
$ws.Range("C71").Value = "[name=""'Captain'""]   We're a long way from the nearest city, and we barely have any signal here. Even though we probably can't count on backup, the other guys are in the same situation. This is the best chance we're gonna get!`n"
$ws.Range("C73").Value = "[name=""'Captain'""]   What are you scared of? She might know how to fight, but don't forget that we know the terrain better than she does. Just wait 'til she discovers all the traps we laid for her!`n"
$ws.Range("C74").Value = "[name=""'Captain'""]   Treasures exist for us to take! Monsters exist for us to hunt!`n"
$ws.Range("C75").Value = "[name=""'Captain'""]   I will have it all!`n"
